$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Mark two more rows as "OK" in column D (fixes the "page loading" problem tracking)
$ws.Range("D3").Value = "OK"
$ws.Range("D11").Value = "OK"

# Widen column B by 7 characters so the longer text fits
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth + 7

# Move the active selection to the newly filled-in cell
$ws.Range("D11").Select()
